$d = $word.ActiveDocument

# --- Numbers / String section -------------------------------------------
# "Variable created with the -S symbol" -> "Variable created with the S-- symbol"
$r1 = $d.Content.Find.Execute(
    "Variable created with the -S symbol", $true, $false, $false, $false, $false,
    $true, 1, $false, "Variable created with the S-- symbol", 2)
Write-Output "Step1 (S-- symbol): $r1"

# --- array section --------------------------------------------------------
# "Created using the -A identifier" ->
# "Created using the A-- identifier    Note: arrays in SB are dynamic, similar to vectors in C++"
$r2 = $d.Content.Find.Execute(
    "Created using the -A identifier", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Created using the A-- identifier    Note: arrays in SB are dynamic, similar to vectors in C++",
    2)
Write-Output "Step2 (A-- identifier + note): $r2"

# "Ex: -A arry1#" -> "Ex: A-- arry1#   : creates an array with no currently assigned values"
$r3 = $d.Content.Find.Execute(
    "Ex: -A arry1#", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Ex: A-- arry1#   : creates an array with no currently assigned values",
    2)
Write-Output "Step3 (Ex: A-- arry1#): $r3"

# "      -A arry2 {(1, 2, 3)}    [en dash]   1D array[tab]-A arry3{(1, 2, 3) : (4, 5, 6)}   - 2D array"
# -> "      A-- arry2 {(1, 2, 3)}    [en dash]   1D array[tab]A-- arry3{(1, 2, 3) : (4, 5, 6)}   - 2D array"
$old4 = "      -A arry2 {(1, 2, 3)}    " + [char]0x2013 + "   1D array" + [char]9 + `
        "-A arry3{(1, 2, 3) : (4, 5, 6)}   - 2D array"
$new4 = "      A-- arry2 {(1, 2, 3)}    " + [char]0x2013 + "   1D array" + [char]9 + `
        "A-- arry3{(1, 2, 3) : (4, 5, 6)}   - 2D array"
$r4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Output "Step4 (arry2/arry3 A--): $r4"

# Move the "_GoBack" bookmark from after "... shows no decimal amount" down to sit
# between the second "A--" and " arry3{(1, 2, 3) : (4, 5, 6)}   - 2D array" text,
# matching where the cursor was left after the edit above.
$rng = $d.Content
$found5 = $rng.Find.Execute("A-- arry3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $bmStart = $rng.Start + 3
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    Write-Output "Step5 (bookmark moved to $bmStart): True"
} else {
    Write-Output "Step5 (bookmark moved): False"
}

# --- For loop section -------------------------------------------------
# Merge the split runs "3" + "#" + " /c Will increment by 3c\" into a single run.
$old6 = "3# /c Will increment by 3c\"
$r6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2)
Write-Output "Step6 (merge 3# run): $r6"

Write-Output "All edits applied"
